# Applies the "Added data and templates" edit to the Mid-1 hall-ticket
# template: tidies a couple of grammar-checker-split runs, bumps the
# semester from III to II, refreshes the exam schedule table with the
# new September dates/subjects (inserting an extra row for the
# Numerical Ability paper), and drops the stray blank line that used
# to sit right under the table.

$d = $word.ActiveDocument

function Replace-Text($range, [string]$oldText, [string]$newText) {
    $range.Find.Execute($oldText, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $newText, 2) | Out-Null
}

# --- Header block -----------------------------------------------------
# "Approved " / "By" / " " (split by the grammar checker) -> one run.
Replace-Text $d.Content "Approved By " "Approved By "

# "...Mid-1 Examinations, III B.Tech I Sem..." -> "...II B.Tech I Sem..."
Replace-Text $d.Content `
    "Hall Ticket: Mid-1 Examinations, III B.Tech I Sem" `
    "Hall Ticket: Mid-1 Examinations, II B.Tech I Sem"

# "{{PIN" / "}}" (split by the grammar checker) -> one "{{PIN}}" run.
Replace-Text $d.Content "{{PIN}}" "{{PIN}}"

# --- Exam schedule table -----------------------------------------------
$t = $d.Tables.Item(1)

# Row 2: 18-08-2025, Information Retrieval Systems
Replace-Text $t.Cell(2,1).Range "18-08-2025" "01-09-2025"
Replace-Text $t.Cell(2,2).Range "Information Retrieval Systems" "Artificial Intelligence"

# Row 3: 19-08-2025, Computer Networks
Replace-Text $t.Cell(3,1).Range "19-08-2025" "02-09-2025"
Replace-Text $t.Cell(3,2).Range "Computer Networks" "Object Oriented Programming through JAVA"

# Row 4: 20-08-2025, Automata Theory and Compiler Design
Replace-Text $t.Cell(4,1).Range "20-08-2025" "03-09-2025"
Replace-Text $t.Cell(4,2).Range "Automata Theory and Compiler Design" "Database Management Systems"

# Row 5: 21-08-2025, Block Chain Technologies and its Applications
Replace-Text $t.Cell(5,1).Range "21-08-2025" "04-09-2025"
Replace-Text $t.Cell(5,2).Range "Block Chain Technologies and its Applications" "Operating Systems"

# New row inserted before the old "22-08-2025" row: 08-09-2025, Numerical
# Ability and Cognitive Intelligence.
$newRow = $t.Rows.Add($t.Rows.Item(6))
$t.Cell(6,1).Range.Text = "08-09-2025"
$t.Cell(6,2).Range.Text = "Numerical Ability and Cognitive Intelligence"

# Old row 6 (now row 7): 22-08-2025, Internet of Things and Applications
Replace-Text $t.Cell(7,1).Range "22-08-2025" "09-09-2025"
Replace-Text $t.Cell(7,2).Range "Internet of Things and Applications" "Basics of Business Management"

# --- Drop the stray blank "<w:br/>" paragraph right under the table ---
$tEnd = $t.Range.End
$d.Range($tEnd, $tEnd + 2).Delete()

# --- Styles: Default Paragraph Font is no longer semi-hidden ----------
$dpf = $d.Styles.Item("Default Paragraph Font")
$dpf.UnhideWhenUsed = $true
